# Insert a new data row at row 20 (pushing the existing rows 20-46 down to
# 21-47), then populate the new row with its values. This mirrors the
# author's edit: a new weekly "Rabanito" price observation was inserted into
# the middle of the historical series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 20:46 down by one row, creating a blank row 20.
$ws.Rows("20:20").Insert()

# Fill in the newly inserted row 20 with the new observation's data.
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = 44494
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 300000001
$ws.Range("G20").Value = "Rabanito"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = 6000
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = 6000
$ws.Range("N20").Value = "$/docena de paquetes"
$ws.Range("O20").Value = "Provincia de Cautín"
$ws.Range("P20").Value = 500
$ws.Range("Q20").Value = 12
$ws.Range("R20").Value = "Hortaliza"
